$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
for ($r = 11; $r -le 43; $r++) {
    $name = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 16).Value = ".fa-" + $name
}
